# "actualice lista de precios 23-1-26"
#
# GRUPOS_CLIENTES.xlsx / Hoja1
#   - B4 (group "PREMIUM TOP") changes from the bare number 949 to the
#     text/list "949.1218.1959.5625" (same pattern as every other row in
#     column B, which already stores dot-separated client-id lists as
#     text). Giving it a Text number format keeps the leading digits from
#     being re-interpreted as a number and keeps the left-aligned look the
#     cell already had.
#   - The sheet view was scrolled/zoomed while reviewing the update
#     (zoom 70% -> 115%, selection moved from B4 to A4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- data change -----------------------------------------------------
$cell = $ws.Range("B4")
$cell.Value = "949.1218.1959.5625"
$cell.NumberFormat = "@"

# --- view change -------------------------------------------------------
$ws.Activate()
$null = $ws.Range("A4").Select()
$excel.ActiveWindow.Zoom = 115
